$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update data rows (A2:H49) with newly scraped consorcio listings ---
# Force column E ("Total de Parcelas") to Text format so purely numeric values
# (e.g. "138") are stored as text, matching the source data type.
$ws.Range("E2:E49").NumberFormat = "@"

$data = New-Object 'object[,]' 48,8
$data[0,0] = 'CS119'
$data[0,1] = 'Imóveis'
$data[0,2] = 'R$ 69.600,00'
$data[0,3] = 'R$ 22.480,00'
$data[0,4] = '138'
$data[0,5] = 'Itaú'
$data[0,6] = 'Disponível'
$data[0,7] = '138x R$ 479,00'
$data[1,0] = 'CS120'
$data[1,1] = 'Imóveis'
$data[1,2] = 'R$ 70.200,00'
$data[1,3] = 'R$ 44.510,00'
$data[1,4] = '192'
$data[1,5] = 'Porto Seguro'
$data[1,6] = 'Disponível'
$data[1,7] = '192x R$ 399,00'
$data[2,0] = 'CS121'
$data[2,1] = 'Imóveis'
$data[2,2] = 'R$ 110.000,00'
$data[2,3] = 'R$ 70.500,00'
$data[2,4] = '177'
$data[2,5] = 'Porto Seguro'
$data[2,6] = 'Disponível'
$data[2,7] = '177x R$ 469,00'
$data[3,0] = 'CS122'
$data[3,1] = 'Imóveis'
$data[3,2] = 'R$ 111.000,00'
$data[3,3] = 'R$ 67.550,00'
$data[3,4] = '73'
$data[3,5] = 'Itaú'
$data[3,6] = 'Disponível'
$data[3,7] = '73x R$ 1.387,00'
$data[4,0] = 'CS123'
$data[4,1] = 'Imóveis'
$data[4,2] = 'R$ 120.000,00'
$data[4,3] = 'R$ 75.000,00'
$data[4,4] = '56'
$data[4,5] = 'Itaú'
$data[4,6] = 'Disponível'
$data[4,7] = '56x R$ 1.608,00'
$data[5,0] = 'CS124'
$data[5,1] = 'Imóveis'
$data[5,2] = 'R$ 128.000,00'
$data[5,3] = 'R$ 64.400,00'
$data[5,4] = '52'
$data[5,5] = 'Itaú'
$data[5,6] = 'Disponível'
$data[5,7] = '52x R$ 2.333,00'
$data[6,0] = 'CS125'
$data[6,1] = 'Imóveis'
$data[6,2] = 'R$ 128.000,00'
$data[6,3] = 'R$ 64.400,00'
$data[6,4] = '52'
$data[6,5] = 'Itaú'
$data[6,6] = 'Disponível'
$data[6,7] = '52x R$ 2.342,00'
$data[7,0] = 'CS126'
$data[7,1] = 'Imóveis'
$data[7,2] = 'R$ 132.000,00'
$data[7,3] = 'R$ 31.600,00'
$data[7,4] = '169'
$data[7,5] = 'BP Consórcio'
$data[7,6] = 'Disponível'
$data[7,7] = '169x R$ 1.027,00'
$data[8,0] = 'CS127'
$data[8,1] = 'Imóveis'
$data[8,2] = 'R$ 135.000,00'
$data[8,3] = 'R$ 72.750,00'
$data[8,4] = '92'
$data[8,5] = 'Porto Seguro'
$data[8,6] = 'Disponível'
$data[8,7] = '92x R$ 1.395,00'
$data[9,0] = 'CS128'
$data[9,1] = 'Imóveis'
$data[9,2] = 'R$ 137.000,00'
$data[9,3] = 'R$ 73.850,00'
$data[9,4] = '52'
$data[9,5] = 'Itaú'
$data[9,6] = 'Disponível'
$data[9,7] = '52x R$ 2.339,00'
$data[10,0] = 'CS129'
$data[10,1] = 'Imóveis'
$data[10,2] = 'R$ 149.300,00'
$data[10,3] = 'R$ 69.465,00'
$data[10,4] = '136'
$data[10,5] = 'Bradesco'
$data[10,6] = 'Disponível'
$data[10,7] = '136x R$ 854,00'
$data[11,0] = 'CS130'
$data[11,1] = 'Imóveis'
$data[11,2] = 'R$ 150.000,00'
$data[11,3] = 'R$ 75.500,00'
$data[11,4] = '52'
$data[11,5] = 'Itaú'
$data[11,6] = 'Disponível'
$data[11,7] = '52x R$ 2.951,00'
$data[12,0] = 'CS131'
$data[12,1] = 'Imóveis'
$data[12,2] = 'R$ 153.000,00'
$data[12,3] = 'R$ 83.650,00'
$data[12,4] = '83'
$data[12,5] = 'Porto Seguro'
$data[12,6] = 'Disponível'
$data[12,7] = '83x R$ 1.845,00'
$data[13,0] = 'CS132'
$data[13,1] = 'Imóveis'
$data[13,2] = 'R$ 158.000,00'
$data[13,3] = 'R$ 82.900,00'
$data[13,4] = '52'
$data[13,5] = 'Itaú'
$data[13,6] = 'Disponível'
$data[13,7] = '52x R$ 2.894,00'
$data[14,0] = 'CS133'
$data[14,1] = 'Imóveis'
$data[14,2] = 'R$ 187.200,00'
$data[14,3] = 'R$ 104.360,00'
$data[14,4] = '124'
$data[14,5] = 'Santander'
$data[14,6] = 'Disponível'
$data[14,7] = '124x R$ 1.064,00'
$data[15,0] = 'CS134'
$data[15,1] = 'Imóveis'
$data[15,2] = 'R$ 212.000,00'
$data[15,3] = 'R$ 131.600,00'
$data[15,4] = '192'
$data[15,5] = 'Porto Seguro'
$data[15,6] = 'Disponível'
$data[15,7] = '192x R$ 1.342,00'
$data[16,0] = 'CS135'
$data[16,1] = 'Imóveis'
$data[16,2] = 'R$ 219.000,00'
$data[16,3] = 'R$ 140.950,00'
$data[16,4] = '177'
$data[16,5] = 'Porto Seguro'
$data[16,6] = 'Disponível'
$data[16,7] = '177x R$ 938,00'
$data[17,0] = 'CS136'
$data[17,1] = 'Imóveis'
$data[17,2] = 'R$ 474.300,00'
$data[17,3] = 'R$ 258.715,00'
$data[17,4] = '105'
$data[17,5] = 'Porto Seguro'
$data[17,6] = 'Disponível'
$data[17,7] = '105x R$ 3.603,00'
$data[18,0] = 'CS137'
$data[18,1] = 'Imóveis'
$data[18,2] = 'R$ 527.000,00'
$data[18,3] = 'R$ 255.350,00'
$data[18,4] = '180'
$data[18,5] = 'Porto Seguro'
$data[18,6] = 'Disponível'
$data[18,7] = '180x R$ 3.430,00'
$data[19,0] = 'CS138'
$data[19,1] = 'Imóveis'
$data[19,2] = 'R$ 546.000,00'
$data[19,3] = 'R$ 286.300,00'
$data[19,4] = '144'
$data[19,5] = 'Porto Seguro'
$data[19,6] = 'Disponível'
$data[19,7] = '144x R$ 4.370,00'
$data[20,0] = 'CS139'
$data[20,1] = 'Imóveis'
$data[20,2] = 'R$ 1.334.000,00'
$data[20,3] = 'R$ 665.700,00'
$data[20,4] = '144'
$data[20,5] = 'Itaú'
$data[20,6] = 'Disponível'
$data[20,7] = '144x R$ 10.952,00'
$data[21,0] = 'CS140'
$data[21,1] = 'Veículos'
$data[21,2] = 'R$ 6.800,00'
$data[21,3] = 'R$ 2.840,00'
$data[21,4] = '58'
$data[21,5] = 'Itaú'
$data[21,6] = 'Disponível'
$data[21,7] = '58x R$ 90,00'
$data[22,0] = 'CS141'
$data[22,1] = 'Veículos'
$data[22,2] = 'R$ 10.000,00'
$data[22,3] = 'R$ 6.000,00'
$data[22,4] = '33'
$data[22,5] = 'Santander'
$data[22,6] = 'Disponível'
$data[22,7] = '33x R$ 135,00'
$data[23,0] = 'CS142'
$data[23,1] = 'Veículos'
$data[23,2] = 'R$ 18.600,00'
$data[23,3] = 'R$ 11.930,00'
$data[23,4] = '53'
$data[23,5] = 'Itaú'
$data[23,6] = 'Disponível'
$data[23,7] = '53x R$ 365,00'
$data[24,0] = 'CS143'
$data[24,1] = 'Veículos'
$data[24,2] = 'R$ 31.250,00'
$data[24,3] = 'R$ 13.562,50'
$data[24,4] = '43'
$data[24,5] = 'Santander'
$data[24,6] = 'Disponível'
$data[24,7] = '43x R$ 634,00'
$data[25,0] = 'CS144'
$data[25,1] = 'Veículos'
$data[25,2] = 'R$ 35.600,00'
$data[25,3] = 'R$ 19.780,00'
$data[25,4] = '44'
$data[25,5] = 'Itaú'
$data[25,6] = 'Disponível'
$data[25,7] = '44x R$ 709,00'
$data[26,0] = 'CS145'
$data[26,1] = 'Veículos'
$data[26,2] = 'R$ 40.700,00'
$data[26,3] = 'R$ 23.035,00'
$data[26,4] = '49'
$data[26,5] = 'Itaú'
$data[26,6] = 'Disponível'
$data[26,7] = '49x R$ 847,00'
$data[27,0] = 'CS146'
$data[27,1] = 'Veículos'
$data[27,2] = 'R$ 41.000,00'
$data[27,3] = 'R$ 24.050,00'
$data[27,4] = '49'
$data[27,5] = 'Itaú'
$data[27,6] = 'Disponível'
$data[27,7] = '49x R$ 825,00'
$data[28,0] = 'CS147'
$data[28,1] = 'Veículos'
$data[28,2] = 'R$ 67.700,00'
$data[28,3] = 'R$ 35.385,00'
$data[28,4] = '54'
$data[28,5] = 'Itaú'
$data[28,6] = 'Disponível'
$data[28,7] = '54x R$ 1.250,00'
$data[29,0] = 'CS148'
$data[29,1] = 'Veículos'
$data[29,2] = 'R$ 70.900,00'
$data[29,3] = 'R$ 26.545,00'
$data[29,4] = '37'
$data[29,5] = 'Itaú'
$data[29,6] = 'Disponível'
$data[29,7] = '37x R$ 2.424,00'
$data[30,0] = 'CS149'
$data[30,1] = 'Veículos'
$data[30,2] = 'R$ 106.300,00'
$data[30,3] = 'R$ 38.315,00'
$data[30,4] = '37'
$data[30,5] = 'Itaú'
$data[30,6] = 'Disponível'
$data[30,7] = '37x R$ 3.636,00'
$data[31,0] = 'CS150'
$data[31,1] = 'Veículos'
$data[31,2] = 'R$ 109.000,00'
$data[31,3] = 'R$ 67.450,00'
$data[31,4] = '49'
$data[31,5] = 'Porto Seguro'
$data[31,6] = 'Disponível'
$data[31,7] = '49x R$ 1.657,00'
$data[32,0] = 'CS151'
$data[32,1] = 'Veículos'
$data[32,2] = 'R$ 118.600,00'
$data[32,3] = 'R$ 69.930,00'
$data[32,4] = '34'
$data[32,5] = 'Unicoob (Sicoob)'
$data[32,6] = 'Disponível'
$data[32,7] = '34x R$ 2.018,00'
$data[33,0] = 'CS152'
$data[33,1] = 'Veículos'
$data[33,2] = 'R$ 119.000,00'
$data[33,3] = 'R$ 67.950,00'
$data[33,4] = '48'
$data[33,5] = 'Porto Seguro'
$data[33,6] = 'Disponível'
$data[33,7] = '48x R$ 2.250,00'
$data[34,0] = 'CS153'
$data[34,1] = 'Veículos'
$data[34,2] = 'R$ 121.300,00'
$data[34,3] = 'R$ 41.065,00'
$data[34,4] = '29'
$data[34,5] = 'Itaú'
$data[34,6] = 'Disponível'
$data[34,7] = '29x R$ 4.932,00'
$data[35,0] = 'CS154'
$data[35,1] = 'Veículos'
$data[35,2] = 'R$ 125.000,00'
$data[35,3] = 'R$ 68.250,00'
$data[35,4] = '64'
$data[35,5] = 'Porto Seguro'
$data[35,6] = 'Disponível'
$data[35,7] = '64x R$ 1.877,00'
$data[36,0] = 'CS155'
$data[36,1] = 'Veículos'
$data[36,2] = 'R$ 129.000,00'
$data[36,3] = 'R$ 78.450,00'
$data[36,4] = '64'
$data[36,5] = 'Porto Seguro'
$data[36,6] = 'Disponível'
$data[36,7] = '64x R$ 1.665,00'
$data[37,0] = 'CS156'
$data[37,1] = 'Veículos'
$data[37,2] = 'R$ 140.000,00'
$data[37,3] = 'R$ 65.000,00'
$data[37,4] = '41'
$data[37,5] = 'Itaú'
$data[37,6] = 'Disponível'
$data[37,7] = '41x R$ 3.678,00'
$data[38,0] = 'CS157'
$data[38,1] = 'Veículos'
$data[38,2] = 'R$ 201.000,00'
$data[38,3] = 'R$ 105.050,00'
$data[38,4] = '66'
$data[38,5] = 'Porto Seguro'
$data[38,6] = 'Disponível'
$data[38,7] = '66x R$ 2.902,00'
$data[39,0] = 'CS158'
$data[39,1] = 'Veículos'
$data[39,2] = 'R$ 207.000,00'
$data[39,3] = 'R$ 105.350,00'
$data[39,4] = '64'
$data[39,5] = 'Porto Seguro'
$data[39,6] = 'Disponível'
$data[39,7] = '64x R$ 2.966,00'
$data[40,0] = 'CS159'
$data[40,1] = 'Veículos'
$data[40,2] = 'R$ 209.000,00'
$data[40,3] = 'R$ 99.450,00'
$data[40,4] = '57'
$data[40,5] = 'Itaú'
$data[40,6] = 'Disponível'
$data[40,7] = '57x R$ 3.590,00'
$data[41,0] = 'CS160'
$data[41,1] = 'Veículos'
$data[41,2] = 'R$ 210.000,00'
$data[41,3] = 'R$ 99.500,00'
$data[41,4] = '57'
$data[41,5] = 'Itaú'
$data[41,6] = 'Disponível'
$data[41,7] = '57x R$ 3.630,00'
$data[42,0] = 'CS161'
$data[42,1] = 'Veículos'
$data[42,2] = 'R$ 254.000,00'
$data[42,3] = 'R$ 144.700,00'
$data[42,4] = '64'
$data[42,5] = 'Porto Seguro'
$data[42,6] = 'Disponível'
$data[42,7] = '64x R$ 3.520,00'
$data[43,0] = 'CS162'
$data[43,1] = 'Veículos'
$data[43,2] = 'R$ 310.000,00'
$data[43,3] = 'R$ 130.500,00'
$data[43,4] = '47'
$data[43,5] = 'Itaú'
$data[43,6] = 'Disponível'
$data[43,7] = '47x R$ 6.190,00'
$data[44,0] = 'CS163'
$data[44,1] = 'Veículos'
$data[44,2] = 'R$ 312.000,00'
$data[44,3] = 'R$ 132.600,00'
$data[44,4] = '47'
$data[44,5] = 'Itaú'
$data[44,6] = 'Disponível'
$data[44,7] = '47x R$ 6.079,00'
$data[45,0] = 'CS164'
$data[45,1] = 'Veículos'
$data[45,2] = 'R$ 332.200,00'
$data[45,3] = 'R$ 171.610,00'
$data[45,4] = '64'
$data[45,5] = 'Porto Seguro'
$data[45,6] = 'Disponível'
$data[45,7] = '64x R$ 4.841,00'
$data[46,0] = 'CS165'
$data[46,1] = 'Veículos'
$data[46,2] = 'R$ 471.000,00'
$data[46,3] = 'R$ 263.550,00'
$data[46,4] = '48'
$data[46,5] = 'Itaú'
$data[46,6] = 'Disponível'
$data[46,7] = '48x R$ 6.820,00'
$data[47,0] = 'CS166'
$data[47,1] = 'Veículos'
$data[47,2] = 'R$ 620.000,00'
$data[47,3] = 'R$ 261.000,00'
$data[47,4] = '47'
$data[47,5] = 'Itaú'
$data[47,6] = 'Disponível'
$data[47,7] = '47x R$ 12.269,00'

$ws.Range("A2:H49").Value = $data

# --- Ensure the newly appended rows keep the trailing empty "Vencimento"/"Observacoes" cells ---
$ws.Range("I37:J49").NumberFormat = "@"
$ws.Range("I37:J49").Value = ""

# --- Widen column F ("Consorcio") from 14 to 18 characters ---
$ws.Columns.Item(6).ColumnWidth = 17.1666666666667

Write-Host "Sheet updated: 48 data rows, column F resized"